$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'309.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-1.83%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'37.97"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-3.63%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.067"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-1.16%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07769"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-4.94%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.357"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-0.20%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.900"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-3.79%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'8.191"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-1.57%"
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9203"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-1.74%"
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1252"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-3.90%"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1881"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-4.51%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.08813"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-2.40%"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03412"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-2.38%"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09713"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.19%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001365"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-2.95%"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006068"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.01%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.573"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-2.07%"
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "'3.094"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-6.25%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'-2.26%"
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").Value = "'0.1279"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-2.82%"
$ws.Range("E20").Style = "Normal"
$ws.Range("B21").Value = "MCDex"
$ws.Range("C21").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D21").Value = "'5.029"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'1.40%"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'4.02%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.02105"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'5,592.37%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.04395"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.63%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'-2.59%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004255"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-10.68%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0001351"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-65.31%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.02131"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-3.71%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04988"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-3.72%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.008056"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'3.82%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.01003"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-3.37%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.1343"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-4.07%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002062"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-1.93%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.008706"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-6.20%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006461"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-6.73%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.02%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.003395"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'17.78%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'-0.13%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.02%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.02%"
$ws.Range("E51").Style = "Normal"
